$d = $word.ActiveDocument

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0137C50D" w14:textId="4A6C53B2" w:rsidR="00770F96" w:rsidRPr="002E712E" w:rsidRDefault="00770F96"><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Raggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> stirbt irgendwann in der Matrix. Dadurch, dass er stirbt, wird er einfach aus der Matrix ausgeklinkt, erhält jedoch nicht seine Erinnerungen zurück und so hilft er den Menschen in der Matrix, immer nach einer Möglichkeit suchend, sich wieder einzuklinken. Der Experimentator</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> versucht irgendwann die Matrix abzuschalten. Das klappt nicht, also versucht er sie zu zerstören. Das funktioniert auch nicht und so klinkt der Experimentator selber sich ein, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Raggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> kommt mit, ohne Wissen des Experimentators. Endschlacht ist ein Kampf gegen ein wildes Sammelsurium aus Zeiten mit dem Experimentator als </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Endboss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($xml4)

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6540FE70" w14:textId="21ACE772" w:rsidR="00920C95" w:rsidRDefault="00920C95"><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">NPCs: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>spawnen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> einfach, ebenfalls oh</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">ne Erinnerung. Irgendwann findet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Raggi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> einen NPC </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Spawner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>. Er findet es ekelhaft.</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> Das Lagerfeuer am Anfang ist ein solcher </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Spawner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($xml3)

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="01319F59" w14:textId="49CD54F4" w:rsidR="005163E9" w:rsidRDefault="005163E9"><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Panik am Anfang lässt sich </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>durch zunehmende</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> Schwärze auf dem Bildschirm darstellen. Ist der Bildschirm schwarz, ist das Spiel vorbei.</w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($xml2)

$xml0 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6B9AE65E" w14:textId="77777777" w:rsidR="000170A2" w:rsidRDefault="002E712E"><w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r w:rsidRPr="002E712E"><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Notizen für Geschichte für Lands </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>f</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> Oblivion</w:t></w:r></w:p>
'@
$d.Paragraphs(1).Range.InsertXML($xml0)
